{"js": "// Update the regression-output numbers that changed between the two\n// builds of the post (scikit-learn section x2, custom-OLS section,\n// and the \"from scratch\" array([...]) section).\nconst replacements = [\n  [\"1.0038 and coefficients\", \"0.9311 and coefficients\"],\n  [\"2.9441, and\", \"3.0374, and\"],\n  [\"1.9985\", \"1.9636\"],\n  [\"1.0032 and coefficients\", \"0.931 and coefficients\"],\n  [\"2.9434, and\", \"3.0363, and\"],\n  [\"1.9978\", \"1.9632\"],\n  [\"array([0.98754779]) and coefficients\", \"array([0.91697402]) and coefficients\"],\n  [\"array([2.95617679]), and\", \"array([3.04908545]), and\"],\n  [\"array([2.00473963])\", \"array([1.98015568])\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the regression-output numbers that changed between the two\n# builds of the post (scikit-learn section x2, custom-OLS section,\n# and the \"from scratch\" array([...]) section).\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"1.0038 and coefficients\", \"0.9311 and coefficients\"),\n    @(\"2.9441, and\", \"3.0374, and\"),\n    @(\"1.9985\", \"1.9636\"),\n    @(\"1.0032 and coefficients\", \"0.931 and coefficients\"),\n    @(\"2.9434, and\", \"3.0363, and\"),\n    @(\"1.9978\", \"1.9632\"),\n    @(\"array([0.98754779]) and coefficients\", \"array([0.91697402]) and coefficients\"),\n    @(\"array([2.95617679]), and\", \"array([3.04908545]), and\"),\n    @(\"array([2.00473963])\", \"array([1.98015568])\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, $wdReplaceAll)\n}\n"}
